$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(15, 8).Value = 356080.94
$ws.Cells.Item(15, 9).Value = 356080.94
$ws.Cells.Item(15, 11).Value = 1068242.82
$ws.Cells.Item(15, 13).Value = -1068073.82
$ws.Cells.Item(43, 8).Value = 2095.4
$ws.Cells.Item(43, 9).Value = 2850
$ws.Cells.Item(43, 10).Value = 1592.3334
$ws.Cells.Item(43, 11).Value = 2850
$ws.Cells.Item(43, 12).Value = 1592.3334
$ws.Cells.Item(43, 13).Value = -2781
$ws.Cells.Item(43, 14).Value = -1730.3334
$ws.Cells.Item(58, 8).Value = 1005.6875
$ws.Cells.Item(58, 9).Value = 330.07693
$ws.Cells.Item(58, 10).Value = 3933.3333
$ws.Cells.Item(58, 11).Value = 990.2307900000001
$ws.Cells.Item(58, 12).Value = 11799.9999
$ws.Cells.Item(58, 13).Value = -840.2307900000001
$ws.Cells.Item(58, 14).Value = -12099.9999
$ws.Cells.Item(132, 8).Value = 19871
$ws.Cells.Item(132, 9).Value = 2927.5
$ws.Cells.Item(132, 10).Value = 31166.666
$ws.Cells.Item(132, 11).Value = 8782.5
$ws.Cells.Item(132, 12).Value = 93499.99800000001
$ws.Cells.Item(132, 13).Value = -6252.5
$ws.Cells.Item(132, 14).Value = -98559.99800000001
$ws.Cells.Item(135, 8).Value = 497.45456
$ws.Cells.Item(135, 9).Value = 184.94118
$ws.Cells.Item(135, 10).Value = 1560
$ws.Cells.Item(135, 11).Value = 1664.47062
$ws.Cells.Item(135, 12).Value = 14040
$ws.Cells.Item(135, 13).Value = 870.5293799999999
$ws.Cells.Item(135, 14).Value = -19110
$ws.Cells.Item(137, 8).Value = 230543.98
$ws.Cells.Item(137, 9).Value = 8698.526
$ws.Cells.Item(137, 10).Value = 338622.53
$ws.Cells.Item(137, 11).Value = 26095.578
$ws.Cells.Item(137, 12).Value = 1015867.59
$ws.Cells.Item(137, 13).Value = -23545.578
$ws.Cells.Item(137, 14).Value = -1020967.59

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(3, 8).Value = 500
$ws.Cells.Item(3, 9).Value = 500
$ws.Cells.Item(3, 11).Value = 500
$ws.Cells.Item(3, 13).Value = -385
$ws.Cells.Item(32, 8).Value = 682.26
$ws.Cells.Item(32, 9).Value = 687.8969
$ws.Cells.Item(32, 10).Value = 500
$ws.Cells.Item(32, 11).Value = 687.8969
$ws.Cells.Item(32, 12).Value = 500
$ws.Cells.Item(32, 13).Value = -400.8969
$ws.Cells.Item(32, 14).Value = -1074
$ws.Cells.Item(45, 8).Value = 899
$ws.Cells.Item(45, 9).Value = 899
$ws.Cells.Item(45, 11).Value = 899
$ws.Cells.Item(45, 13).Value = -522
$ws.Cells.Item(74, 8).Value = 18825864
$ws.Cells.Item(74, 9).Value = 22890024
$ws.Cells.Item(74, 10).Value = 15922892
$ws.Cells.Item(74, 11).Value = 22890024
$ws.Cells.Item(74, 12).Value = 15922892
$ws.Cells.Item(74, 13).Value = -22889150
$ws.Cells.Item(74, 14).Value = -15924640
$ws.Cells.Item(77, 8).Value = 18825864
$ws.Cells.Item(77, 9).Value = 22890024
$ws.Cells.Item(77, 10).Value = 15922892
$ws.Cells.Item(77, 11).Value = 114450120
$ws.Cells.Item(77, 12).Value = 79614460
$ws.Cells.Item(77, 13).Value = -114445752
$ws.Cells.Item(77, 14).Value = -79623196

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(63, 8).Value = 36000
$ws.Cells.Item(63, 10).Value = 36000
$ws.Cells.Item(63, 12).Value = 36000
$ws.Cells.Item(63, 14).Value = -37372
$ws.Cells.Item(66, 8).Value = 36000
$ws.Cells.Item(66, 10).Value = 36000
$ws.Cells.Item(66, 12).Value = 108000
$ws.Cells.Item(66, 14).Value = -114864

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(3, 8).Value = 540.2857
$ws.Cells.Item(3, 9).Value = 100
$ws.Cells.Item(3, 10).Value = 716.4
$ws.Cells.Item(3, 11).Value = 100
$ws.Cells.Item(3, 12).Value = 716.4
$ws.Cells.Item(3, 13).Value = 13
$ws.Cells.Item(3, 14).Value = -942.4
$ws.Cells.Item(31, 8).Value = 2638.5286
$ws.Cells.Item(31, 9).Value = 1081.6207
$ws.Cells.Item(31, 10).Value = 3739.756
$ws.Cells.Item(31, 11).Value = 1081.6207
$ws.Cells.Item(31, 12).Value = 3739.756
$ws.Cells.Item(31, 13).Value = -786.6206999999999
$ws.Cells.Item(31, 14).Value = -4329.755999999999
$ws.Cells.Item(34, 8).Value = 2638.5286
$ws.Cells.Item(34, 9).Value = 1081.6207
$ws.Cells.Item(34, 10).Value = 3739.756
$ws.Cells.Item(34, 11).Value = 1081.6207
$ws.Cells.Item(34, 12).Value = 3739.756
$ws.Cells.Item(34, 13).Value = -879.6206999999999
$ws.Cells.Item(34, 14).Value = -4143.755999999999
$ws.Cells.Item(58, 8).Value = 5135.3145
$ws.Cells.Item(58, 9).Value = 8817.076999999999
$ws.Cells.Item(58, 10).Value = 2959.7273
$ws.Cells.Item(58, 11).Value = 8817.076999999999
$ws.Cells.Item(58, 12).Value = 2959.7273
$ws.Cells.Item(58, 13).Value = -8614.076999999999
$ws.Cells.Item(58, 14).Value = -3365.7273
$ws.Cells.Item(86, 8).Value = 3824.9546
$ws.Cells.Item(86, 9).Value = 2334.3125
$ws.Cells.Item(86, 11).Value = 2334.3125
$ws.Cells.Item(86, 13).Value = -1211.3125
$ws.Cells.Item(89, 8).Value = 3824.9546
$ws.Cells.Item(89, 9).Value = 2334.3125
$ws.Cells.Item(89, 11).Value = 11671.5625
$ws.Cells.Item(89, 13).Value = -6055.5625
$ws.Cells.Item(94, 8).Value = 4547
$ws.Cells.Item(94, 9).Value = 990
$ws.Cells.Item(94, 10).Value = 6206.933
$ws.Cells.Item(94, 11).Value = 990
$ws.Cells.Item(94, 12).Value = 6206.933
$ws.Cells.Item(94, 13).Value = -539
$ws.Cells.Item(94, 14).Value = -7108.933
$ws.Cells.Item(99, 8).Value = 37647.82
$ws.Cells.Item(99, 9).Value = 57362.168
$ws.Cells.Item(99, 10).Value = 2162
$ws.Cells.Item(99, 11).Value = 57362.168
$ws.Cells.Item(99, 12).Value = 2162
$ws.Cells.Item(99, 13).Value = -55864.168
$ws.Cells.Item(99, 14).Value = -5158
$ws.Cells.Item(105, 8).Value = 1198.5454
$ws.Cells.Item(105, 9).Value = 926.75
$ws.Cells.Item(105, 11).Value = 926.75
$ws.Cells.Item(105, 13).Value = 820.25
$ws.Cells.Item(126, 8).Value = 37647.82
$ws.Cells.Item(126, 9).Value = 57362.168
$ws.Cells.Item(126, 10).Value = 2162
$ws.Cells.Item(126, 11).Value = 172086.504
$ws.Cells.Item(126, 12).Value = 6486
$ws.Cells.Item(126, 13).Value = -169616.504
$ws.Cells.Item(126, 14).Value = -11426
$ws.Cells.Item(132, 8).Value = 41669676
$ws.Cells.Item(132, 9).Value = 55557570
$ws.Cells.Item(132, 10).Value = 5999.3335
$ws.Cells.Item(132, 11).Value = 166672710
$ws.Cells.Item(132, 12).Value = 17998.0005
$ws.Cells.Item(132, 13).Value = -166670180
$ws.Cells.Item(132, 14).Value = -23058.0005
$ws.Cells.Item(134, 8).Value = 41668984
$ws.Cells.Item(134, 9).Value = 100001020
$ws.Cells.Item(134, 10).Value = 3244.7144
$ws.Cells.Item(134, 11).Value = 300003060
$ws.Cells.Item(134, 12).Value = 9734.143199999999
$ws.Cells.Item(134, 13).Value = -300000525
$ws.Cells.Item(134, 14).Value = -14804.1432
$ws.Cells.Item(136, 8).Value = 5135.3145
$ws.Cells.Item(136, 9).Value = 8817.076999999999
$ws.Cells.Item(136, 10).Value = 2959.7273
$ws.Cells.Item(136, 11).Value = 26451.231
$ws.Cells.Item(136, 12).Value = 8879.1819
$ws.Cells.Item(136, 13).Value = -23901.231
$ws.Cells.Item(136, 14).Value = -13979.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(13, 8).Value = 600.1667
$ws.Cells.Item(13, 9).Value = 100
$ws.Cells.Item(13, 10).Value = 850.25
$ws.Cells.Item(13, 11).Value = 300
$ws.Cells.Item(13, 12).Value = 2550.75
$ws.Cells.Item(13, 13).Value = -132
$ws.Cells.Item(13, 14).Value = -2886.75
$ws.Cells.Item(23, 8).Value = 506.125
$ws.Cells.Item(23, 9).Value = 50.75
$ws.Cells.Item(23, 10).Value = 597.2
$ws.Cells.Item(23, 11).Value = 152.25
$ws.Cells.Item(23, 12).Value = 1791.6
$ws.Cells.Item(23, 13).Value = 82.75
$ws.Cells.Item(23, 14).Value = -2261.6
$ws.Cells.Item(113, 8).Value = 687.32434
$ws.Cells.Item(113, 9).Value = 668.88
$ws.Cells.Item(113, 10).Value = 725.75
$ws.Cells.Item(113, 11).Value = 2006.64
$ws.Cells.Item(113, 12).Value = 2177.25
$ws.Cells.Item(113, 13).Value = 163.3600000000001
$ws.Cells.Item(113, 14).Value = -6517.25
$ws.Cells.Item(131, 8).Value = 1122.5
$ws.Cells.Item(131, 9).Value = 814.375
$ws.Cells.Item(131, 10).Value = 1239.881
$ws.Cells.Item(131, 11).Value = 2443.125
$ws.Cells.Item(131, 12).Value = 3719.643
$ws.Cells.Item(131, 13).Value = 2596.875
$ws.Cells.Item(131, 14).Value = -13799.643
$ws.Cells.Item(132, 8).Value = 1730.7693
$ws.Cells.Item(132, 9).Value = 2292.5
$ws.Cells.Item(132, 11).Value = 20632.5
$ws.Cells.Item(132, 13).Value = -18102.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(3, 8).Value = 3403667.8
$ws.Cells.Item(3, 9).Value = 10202041
$ws.Cells.Item(3, 10).Value = 4481.2
$ws.Cells.Item(3, 11).Value = 10202041
$ws.Cells.Item(3, 12).Value = 4481.2
$ws.Cells.Item(3, 13).Value = -10201925
$ws.Cells.Item(3, 14).Value = -4713.2
$ws.Cells.Item(132, 8).Value = 3910512.8
$ws.Cells.Item(132, 9).Value = 6948934.5
$ws.Cells.Item(132, 10).Value = 3970.4285
$ws.Cells.Item(132, 11).Value = 20846803.5
$ws.Cells.Item(132, 12).Value = 11911.2855
$ws.Cells.Item(132, 13).Value = -20844273.5
$ws.Cells.Item(132, 14).Value = -16971.2855

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(3, 8).Value = 143638.64
$ws.Cells.Item(3, 9).Value = 182312.81
$ws.Cells.Item(3, 11).Value = 182312.81
$ws.Cells.Item(3, 13).Value = -182198.81
$ws.Cells.Item(11, 8).Value = 70005
$ws.Cells.Item(11, 10).Value = 70005
$ws.Cells.Item(11, 12).Value = 70005
$ws.Cells.Item(11, 14).Value = -70289
$ws.Cells.Item(58, 8).Value = 25000
$ws.Cells.Item(58, 10).Value = 25000
$ws.Cells.Item(58, 12).Value = 25000
$ws.Cells.Item(58, 14).Value = -25616
$ws.Cells.Item(68, 8).Value = 30000
$ws.Cells.Item(68, 10).Value = 30000
$ws.Cells.Item(68, 12).Value = 30000
$ws.Cells.Item(68, 14).Value = -31622
$ws.Cells.Item(71, 8).Value = 30000
$ws.Cells.Item(71, 10).Value = 30000
$ws.Cells.Item(71, 12).Value = 90000
$ws.Cells.Item(71, 14).Value = -98112
$ws.Cells.Item(126, 8).Value = 1203.091
$ws.Cells.Item(126, 9).Value = 1046.5172
$ws.Cells.Item(126, 10).Value = 2338.25
$ws.Cells.Item(126, 11).Value = 3139.5516
$ws.Cells.Item(126, 12).Value = 7014.75
$ws.Cells.Item(126, 13).Value = -669.5515999999998
$ws.Cells.Item(126, 14).Value = -11954.75
$ws.Cells.Item(132, 8).Value = 3231
$ws.Cells.Item(132, 9).Value = 2465.3
$ws.Cells.Item(132, 10).Value = 3927.0908
$ws.Cells.Item(132, 11).Value = 7395.900000000001
$ws.Cells.Item(132, 12).Value = 11781.2724
$ws.Cells.Item(132, 13).Value = -4865.900000000001
$ws.Cells.Item(132, 14).Value = -16841.2724
$ws.Cells.Item(136, 8).Value = 19657850
$ws.Cells.Item(136, 9).Value = 27328710
$ws.Cells.Item(136, 10).Value = 736393.4
$ws.Cells.Item(136, 11).Value = 81986130
$ws.Cells.Item(136, 12).Value = 2209180.2
$ws.Cells.Item(136, 13).Value = -81983580
$ws.Cells.Item(136, 14).Value = -2214280.2
